# ============================================================
# Edit: add "2022-Q1" sheet with fund holdings data, inserted
# right before the "总计" summary sheet, and append a new
# 2022-Q1 row to the "总计" sheet's history table.
# ============================================================

$wb = $excel.ActiveWorkbook

# The "总计" sheet is currently the last sheet (position 6).
$totalName = $wb.Worksheets.Item($wb.Worksheets.Count).Name

# A sheet used as a formatting template for the new quarter sheet
# (any existing quarterly fund-holdings sheet, e.g. "2021-Q4").
$templateSheet = $wb.Worksheets.Item(5)

# --- 1. Insert the new "2022-Q1" sheet right before "总计" ---------------
$q1 = $wb.Worksheets.Add($wb.Worksheets.Item($totalName))
$q1.Name = "2022-Q1"

# Worksheets.Item(<position>) is position-based, so any handle grabbed by
# position before this insert (e.g. a captured $totalSheet) would now
# silently point at the newly-inserted sheet instead. Re-resolve "总计" by
# NAME, after the insert, so every later write lands on the right sheet.
$totalSheet = $wb.Worksheets.Item($totalName)

# --- 2. Header row (basically identical across all quarter sheets) -------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$cols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt 7; $i = $i + 1) {
    $q1.Range($cols[$i] + "1").Value = $headers[$i]
}
$templateSheet.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

# --- 3. Fund holdings data rows -------------------------------------------
# pipe-delimited: 基金代码|基金名称|基金规模|股票总仓位|仓位占比|持有市值(亿元)|仓位排名
$fundRows = @(
    "010336|中欧悦享生活混合A|44.28|90.44|8.85|3.9188|3",
    "040008|华安策略优选混合|52.77|92.07|4.17|2.2005|8",
    "002621|中欧消费主题股票A|19.29|88.29|8.77|1.6917|3",
    "011278|华夏内需驱动混合A|34.32|86.05|3.70|1.2698|4",
    "160314|华夏行业混合(LOF)|21.89|92.15|5.27|1.1536|3",
    "005449|华夏行业龙头混合|11.19|83.96|6.11|0.6837|5",
    "002697|中欧消费主题股票C|6.20|88.29|8.77|0.5437|3",
    "011282|华夏消费龙头混合A|16.37|87.73|2.64|0.4322|10",
    "010020|华夏线上经济主题精选混合|16.45|84.47|2.59|0.4261|9",
    "010852|中欧内需成长混合型证券投资基金A|5.23|91.46|8.12|0.4247|2",
    "005620|中欧品质消费股票A|3.74|90.47|9.03|0.3377|2",
    "420005|天弘周期策略混合|5.25|89.31|5.78|0.3034|4",
    "007202|天弘优质成长企业精选混合|4.81|92.52|4.83|0.2323|7",
    "005621|中欧品质消费股票C|1.11|90.47|9.03|0.1002|2",
    "010337|中欧悦享生活混合C|1.08|90.44|8.85|0.0956|3",
    "010692|华夏核心价值混合A|2.49|87.86|2.70|0.0672|10",
    "010853|中欧内需成长混合型证券投资基金C|0.67|91.46|8.12|0.0544|2",
    "004008|中融鑫思路灵活配置混合A|3.91|35.78|1.35|0.0528|10",
    "004694|天弘策略精选灵活配置混合A|1.11|80.93|3.82|0.0424|5",
    "011279|华夏内需驱动混合C|0.95|86.05|3.70|0.0352|4",
    "011283|华夏消费龙头混合C|1.25|87.73|2.64|0.0330|10",
    "004009|中融鑫思路灵活配置混合C|2.36|35.78|1.35|0.0319|10",
    "010693|华夏核心价值混合C|0.56|87.86|2.70|0.0151|10",
    "004748|天弘策略精选灵活配置混合C|0.08|80.93|3.82|0.0031|5",
    "011494|华泰紫金丰和偏债混合型发起式证券投资基金A|0.15|25.30|0.81|0.0012|6",
    "011495|华泰紫金丰和偏债混合型发起式证券投资基金C|0.03|25.30|0.81|0.0002|6"
)

$r = 2
foreach ($line in $fundRows) {
    $p = $line.Split("|")

    $q1.Range("A" + $r).Value = ($r - 2)

    $q1.Range("B" + $r).NumberFormat = "@"
    $q1.Range("B" + $r).Value = $p[0]

    $q1.Range("C" + $r).Value = $p[1]

    $q1.Range("D" + $r).NumberFormat = "@"
    $q1.Range("D" + $r).Value = $p[2]

    $q1.Range("E" + $r).NumberFormat = "@"
    $q1.Range("E" + $r).Value = $p[3]

    $q1.Range("F" + $r).NumberFormat = "@"
    $q1.Range("F" + $r).Value = $p[4]

    $q1.Range("G" + $r).NumberFormat = "@"
    $q1.Range("G" + $r).Value = $p[5]

    $q1.Range("H" + $r).Value = [int]$p[6]

    $r = $r + 1
}

# Column-A row-index style (centered/bordered, like every other quarter sheet)
$templateSheet.Range("A2:A27").Copy()
$q1.Range("A2:A27").PasteSpecial(-4122)

# --- 4. Update the "总计" (grand total) history sheet ---------------------
# Shift existing rows 2-6 down to rows 3-7 (bottom-up so we never clobber
# a row before it has been read).
for ($row = 6; $row -ge 2; $row = $row - 1) {
    $dest = $row + 1
    $totalSheet.Range("A" + $dest).Value = $totalSheet.Range("A" + $row).Value()
    $totalSheet.Range("B" + $dest).Value = $totalSheet.Range("B" + $row).Value()
    $totalSheet.Range("C" + $dest).Value = $totalSheet.Range("C" + $row).Value()
    $totalSheet.Range("D" + $dest).Value = $totalSheet.Range("D" + $row).Value()
}
# Row 7 is brand new territory (sheet used to stop at row 6) - copy its
# number-column format from row 6 so it matches the rest of the table.
$totalSheet.Range("A6").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)

# New row 2: the 2022-Q1 summary entry.
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 26
$totalSheet.Range("D2").Value = 14.15

# Re-number the index column (A) for the shifted rows: 0..4 -> 1..5.
for ($row = 3; $row -le 7; $row = $row + 1) {
    $totalSheet.Range("A" + $row).Value = $row - 2
}

# --- 5. Restore original active sheet/selection --------------------------
$wb.Worksheets.Item(1).Activate()
